# Update to GAS hierarchy
# - Sheet1!B2 value changes from 20 to 10
# - Active selection moves from A4 to B3 (the last edited/selected cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 10

$ws.Range("B3").Select()
